$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in the new row of data (LD48 / Bouldaouf dash)
$ws.Range("A8").Value = "LD48"
$ws.Range("B8").Value = "Bouldaouf dash"
$ws.Range("C8").Value = 3.235
$ws.Range("D8").Value = 3.103
$ws.Range("E8").Value = 2.765

# F8 and G8 hold numeric-looking text (not real numbers), so build them as
# text via a helper cell/formula and paste the resulting value in, which
# keeps the cell's existing (General) style intact.
$helper = $ws.Range("Z1")

$helper.Formula = "=""2.603"""
$helper.Copy()
$ws.Range("F8").PasteSpecial(-4163)  # xlPasteValues

$helper.Formula = "=""2.735"""
$helper.Copy()
$ws.Range("G8").PasteSpecial(-4163)  # xlPasteValues

$helper.Clear()

$ws.Range("H8").Value = 2.838
$ws.Range("I8").Value = 2.758

# Update the selected cell to B8, matching the saved selection state
$ws.Range("B8").Select()
